$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = -0.6727554909559251

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.7201281033363175
$ws.Range("D5").Value = -0.7501789147931492

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = -0.8147022851187192

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.6255382897778456
$ws.Range("D9").Value = -0.6569424854293864
